# Updated parameters file to clarify the usage of 0 for no vlan tagging scenarios
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datafill")

# Update the Description (column D) text for the N6 VLAN id parameters (DNN1-DNN10)
$ws.Range("D96").Value = "VLAN identifier for traffic on the first N6 (SGi) Data network. Use value as 0 if you don't plan on implementing vlan tagging on the ASE. Do not leave Blank."
$ws.Range("D97").Value = "Optional VLAN identifier for traffic on the second N6 (SGi) Data network. Use value as 0 if you don't plan on implementing vlan tagging on the ASE. If you do not need this DNN/APN, leave the cell empty."
$ws.Range("D98").Value = "Optional VLAN identifier for traffic on the third N6 (SGi) Data network. Use value as 0 if you don't plan on implementing vlan tagging on the ASE. If you do not need this DNN/APN, leave the cell empty."
$ws.Range("D99").Value = "Optional VLAN identifier for traffic on the fourth N6 (SGi) Data network. Use value as 0 if you don't plan on implementing vlan tagging on the ASE. If you do not need this DNN/APN, leave the cell empty."
$ws.Range("D100").Value = "Optional VLAN identifier for traffic on the fifth N6 (SGi) Data network. Use value as 0 if you don't plan on implementing vlan tagging on the ASE. If you do not need this DNN/APN, leave the cell empty."
$ws.Range("D101").Value = "Optional VLAN identifier for traffic on the sixth N6 (SGi) Data network. Use value as 0 if you don't plan on implementing vlan tagging on the ASE. If you do not need this DNN/APN, leave the cell empty."
$ws.Range("D102").Value = "Optional VLAN identifier for traffic on the seventh N6 (SGi) Data network. Use value as 0 if you don't plan on implementing vlan tagging on the ASE. If you do not need this DNN/APN, leave the cell empty."
$ws.Range("D103").Value = "Optional VLAN identifier for traffic on the eighth N6 (SGi) Data network. Use value as 0 if you don't plan on implementing vlan tagging on the ASE. If you do not need this DNN/APN, leave the cell empty."
$ws.Range("D104").Value = "Optional VLAN identifier for traffic on the ninth N6 (SGi) Data network. Use value as 0 if you don't plan on implementing vlan tagging on the ASE. If you do not need this DNN/APN, leave the cell empty."
$ws.Range("D105").Value = "Optional VLAN identifier for traffic on the tenth N6 (SGi) Data network. Use value as 0 if you don't plan on implementing vlan tagging on the ASE. If you do not need this DNN/APN, leave the cell empty."

# Restore the active selection/scroll position on the Datafill sheet after the edit
$ws.Activate()
$win = $excel.ActiveWindow
for ($i = 1; $i -le $win.Panes.Count; $i++) {
  $p = $win.Panes.Item($i)
  $p.ScrollRow = 17
  $p.ScrollColumn = 2
}
$ws.Range("D43").Select()
